$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the "Service locator" list item. In the original document
#    this paragraph holds "Service locator" followed by a manual line
#    break (<w:br/>). We split it into two paragraphs right after the
#    "Service locator" text (i.e. right before the line break), then
#    add a new bullet "Complex data heavy logic (consider DB driven
#    method)" followed by a space ahead of the (now relocated) break.
# ------------------------------------------------------------------

$findRange = $d.Content.Duplicate
$findRange.Find.Execute("Service locator", $false, $false, $false, $false,
                         $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $findRange.End

$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# The new paragraph (still carrying the old line-break run) now begins
# one character past the old split point (the freshly-inserted
# paragraph mark occupies that slot).
$insPos = $splitPos + 1

# ------------------------------------------------------------------
# 2. Insert the trailing space run first (it picks up the paragraph's
#    own formatting - Times New Roman / 13.5pt / black - exactly like
#    the run that used to hold the line break), then prepend the new
#    "Complex data..." run in front of it, explicitly Arial / #222222
#    with the web-style shading seen elsewhere in this document.
# ------------------------------------------------------------------

$spaceXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins = $d.Range($insPos, $insPos)
$ins.InsertXML($spaceXml)

$textXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="222222"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Complex data heavy logic (consider DB driven method)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$ins2 = $d.Range($insPos, $insPos)
$ins2.InsertXML($textXml)

# ------------------------------------------------------------------
# 3. Register the "apple-converted-space" character style (added to
#    styles.xml, based on Default Paragraph Font) that this kind of
#    pasted web content normally brings along with it.
# ------------------------------------------------------------------

$charStyle = $d.Styles.Add("apple-converted-space", 2)
$charStyle.BaseStyle = $d.Styles.Item("DefaultParagraphFont")
